$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new "2020" data column (J), mirroring the existing
# yearly columns D:I. Copy column I's formatting (borders/number format/
# font) into column J first so the new cells match the rest of the table,
# then fill in the actual values.
$ws.Range("I3:I12").Copy() | Out-Null
$ws.Range("J3:J12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("J4").Value = 2020
$ws.Range("J5").Value = 253.27664777870578
$ws.Range("J7").Value = 93.236077839070575
$ws.Range("J8").Value = 160
$ws.Range("J10").Value = 69
$ws.Range("J11").Value = 48.5
$ws.Range("J12").Value = 22.8

# J3, J6 and J9 stay blank (they mirror the blank separator cells in I3/I6/I9).

# The selection moves to the newly added header cell.
$ws.Range("J3").Select()
